# The author's new .ttl was generated from the Google Sheet after the
# "skos:prefLabel" metadata row (row 11) was removed from the vocabulary
# sheet. That removal shifts every following row up by one (rows 12-24 ->
# 11-23), so the sheet shrinks from A1:T24 to A1:T23. In addition, in the
# "Identifier" header row (originally row 19, now row 18) the column-C
# header changed from "pav:version" to "qudt:unit".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "skos:prefLabel" metadata row entirely; this shifts
# all rows below it (old rows 12-24) up by one, which both produces the
# row content shown in the diff and reduces the sheet dimension to T23.
$ws.Rows(11).Delete()

# Update the "Identifier" header row's unit-of-measure column header.
$ws.Range("C18").Value = "qudt:unit"
